# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must stay as text
# (matching the original inlineStr cell type). We briefly force a text
# number format so Excel does not reinterpret them as numeric values,
# then restore the default "Normal" style so no stray formatting remains.
$textCells = @("D5", "D6", "D8", "D12", "D14", "D16", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D46", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.594.00"
$ws.Range("E2").Value = "  +3.43%  "
$ws.Range("D3").Value = "3.290.73"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "574.54"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "176.99"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  +2.97%  "
$ws.Range("D9").Value = "3.285.13"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "45.63"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "701.70"
$ws.Range("E14").Value = "  +11.63%  "
$ws.Range("D15").Value = "3.819.88"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "8.35"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "67.656.29"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "3.286.93"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "17.36"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "10.78"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").Value = "0.891"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").Value = "16.99"
$ws.Range("E23").Value = "  -5.60%  "
$ws.Range("D24").Value = "5.15"
$ws.Range("E24").Value = "  +3.96%  "
$ws.Range("D25").Value = "98.61"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").Value = "2.71"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "32.97"
$ws.Range("E29").Value = "  +7.91%  "
$ws.Range("D30").Value = "8.45"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").Value = "6.70"
$ws.Range("E31").Value = "  +4.09%  "
$ws.Range("D32").Value = "580.60"
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("D33").Value = "3.872.44"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").Value = "10.81"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "0.104"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "55.35"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  -9.35%  "
$ws.Range("D39").Value = "0.128"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").Value = "2.61"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").Value = "3.11"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "3.34"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "31.76"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "0.0₃0675"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("E49").Value = "  +10.52%  "
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "128.07"
$ws.Range("E51").Value = "  -1.27%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
